# Add Sam / Samuel Liska's weekly availability (row 9) and
# update the last-used cell selection, matching the upstream edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 already has "Samuel Liska" in A9; fill in the days Mon-Sun (B9:H9).
# "4pm-MN" is entered before "3pm-MN" so the new shared-string entries land
# in the same order as the source workbook (4pm-MN=29, 3pm-MN=30).
$ws.Range("B9").Value = "2pm-MN"
$ws.Range("D9").Value = "4pm-MN"
$ws.Range("C9").Value = "3pm-MN"
$ws.Range("E9").Value = "3pm-MN"
$ws.Range("F9").Value = "2pm-MN"
$ws.Range("G9").Value = "2pm-MN"
$ws.Range("H9").Value = "2pm-MN"

# Mirror the saved cursor/selection position from the source workbook
$ws.Range("K13").Select()
